# Apply the diff:
#  1. Insert a new worksheet "current_total_expense_base" right after "Estimated"
#     (i.e. as the 2nd sheet), populated with a snapshot of the current "Actual"
#     data plus a new "Cumulative_Quantity" column of zeros.
#  2. Leave "Actual" / "App_Layout" / "actual_cost_v1" / "current_total_expense_v1"
#     / "expense_growth_rate_v1" sheets' data untouched (only their on-disk
#     position/relationship ids shift because of the insertion, which Excel
#     handles automatically).
#  3. Zero out the Quantity/Total/GrandTotal (columns B, D, E) of
#     "planned_estimated_cost_v1" and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- 1. Build the new "current_total_expense_base" sheet -------------------
# Worksheets.Add() inserts immediately before the currently active sheet.
# The workbook opens with "Actual" active (bookViews activeTab="1"), so the
# new sheet lands right after "Estimated" and before "Actual" - exactly where
# the target sheet order needs it.
$base = $wb.Worksheets.Add()
$base.Name = "current_total_expense_base"

$headers = @("Commodity", "Quantity", "Cost", "Total", "GrandTotal", "Cumulative_Quantity")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $base.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Snapshot of the "Actual" sheet's A:E columns (as plain values), with a new
# all-zero Cumulative_Quantity column appended.
$rows = @(
    @("credit", 0, 3000, 0, 8774, 0),
    @("rum", 4, 680, 2720, 7080, 0),
    @("cig", 56, 10, 560, 7080, 0),
    @("netflix", 1, 200, 200, 7080, 0),
    @("amazon", 0, 129, 0, 7080, 0),
    @("onedrive", 0, 145, 0, 7080, 0),
    @("internet", 0, 950, 0, 7080, 0),
    @("Tea", 1, 790, 790, 7080, 0),
    @("Tea Mom", 1, 1160, 1160, 7080, 0),
    @("Juice", 0, 110, 0, 7080, 0),
    @("chicken", 1, 700, 700, 7080, 0),
    @("veg", 1, 270, 270, 7080, 0),
    @("ParkFee", 1, 400, 400, 7080, 0),
    @("Amla and giloy", 0, 900, 0, 7080, 0),
    @("Grocery", 1, 1849, 1849, 7080, 0),
    @("Miscellaneous", 1, 125, 125, 7080, 0),
    @("Cash Withdrawl", 0, 1000, 1000, 7080, 0)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $base.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

$base.Range("F4:F5").Select()

# --- 2. "Actual" sheet: just refresh its selection -------------------------
$actual = $wb.Worksheets.Item("Actual")
$actual.Range("A1:E18").Select()

# --- 3. "planned_estimated_cost_v1": zero out Quantity/Total/GrandTotal ----
$planned = $wb.Worksheets.Item("planned_estimated_cost_v1")
$planned.Range("B2:B17").Value = 0
$planned.Range("D2:D17").Value = 0
$planned.Range("E2:E17").Value = 0
$planned.Range("D2:D17").Select()

# Make it the active sheet/tab, matching the target workbook view.
$planned.Activate()
